# Edit script: swap the "B"/"C" quarter rows within each year block (rows 2-65,
# 4 rows per year: A,B,C,D) across columns A:E, then remove columns F and G
# entirely (the "产销率" and "销售量" per-quarter figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs to swap: for every 4-row year block starting at row 2 (A,B,C,D),
# the 2nd and 3rd rows of the block (the "B" and "C" quarter rows) swap places.
$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56), @(59,60), @(63,64)
)

$cols = @("A", "B", "C", "D", "E")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range($col + $r1)
        $cell2 = $ws.Range($col + $r2)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Remove columns F and G (the "产销率" and "销售量" per-quarter columns); the
# remaining data only spans A:E afterward.
$ws.Range("F1:G65").Delete()

$ws.Range("A1").Select()
